$d = $word.ActiveDocument

# 1. Replace the placeholder ID text (and drop the trailing space run) in the
#    first paragraph: "**ID__AFFARS_pgi_5315_topic_43__ID** " ->
#    "**ID__AFFARS_SMC_PGI_5315_404_1_90__ID**"
$d.Content.Find.Execute("**ID__AFFARS_pgi_5315_topic_43__ID** ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_SMC_PGI_5315_404_1_90__ID**", 2)

# 2. Update paragraph formatting on the first paragraph: indent + paragraph border spacing.
$p = $d.Paragraphs(1)
$pf = $p.Range.ParagraphFormat

$pf.LeftIndent = 11.25  # 225 twips

$borders = $pf.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5
